# Scheduled-runner price refresh: updates market-price-derived columns
# (currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ /
# LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ) across the
# ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR leve-profit sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 672.6585
$ws.Range("J17").Value = 672.6585
$ws.Range("L17").Value = 2017.9755
$ws.Range("N17").Value = -2353.9755

$ws.Range("H38").Value = 34
$ws.Range("I38").Value = 34
$ws.Range("K38").Value = 102
$ws.Range("M38").Value = 270

$ws.Range("H109").Value = 43664.168
$ws.Range("J109").Value = 43664.168
$ws.Range("L109").Value = 43664.168
$ws.Range("N109").Value = -46438.168

$ws.Range("H110").Value = 61885.668
$ws.Range("J110").Value = 61885.668
$ws.Range("L110").Value = 61885.668
$ws.Range("N110").Value = -70065.66800000001

$ws.Range("H134").Value = 98985.336
$ws.Range("J134").Value = 98985.336
$ws.Range("L134").Value = 98985.336
$ws.Range("N134").Value = -109125.336

$ws.Range("H136").Value = 77977.336
$ws.Range("J136").Value = 77977.336
$ws.Range("L136").Value = 77977.336
$ws.Range("N136").Value = -88177.336

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H104").Value = 33659.25
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 33659.25
$ws.Range("K104").Value = 0
$ws.Range("L104").Value = 33659.25
$ws.Range("M104").ClearContents()
$ws.Range("N104").Value = -40647.25

$ws.Range("H107").Value = 58129.2
$ws.Range("J107").Value = 58129.2
$ws.Range("L107").Value = 58129.2
$ws.Range("N107").Value = -65809.2

$ws.Range("H108").Value = 73666
$ws.Range("J108").Value = 87332
$ws.Range("L108").Value = 87332
$ws.Range("N108").Value = -95012

$ws.Range("H118").Value = 54996.8
$ws.Range("J118").Value = 54996.8
$ws.Range("L118").Value = 54996.8
$ws.Range("N118").Value = -58310.8

$ws.Range("H121").Value = 87108.14
$ws.Range("J121").Value = 87108.14
$ws.Range("L121").Value = 87108.14
$ws.Range("N121").Value = -90602.14

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H6").Value = 54540
$ws.Range("J6").Value = 54540
$ws.Range("L6").Value = 54540
$ws.Range("N6").Value = -54766

$ws.Range("H55").Value = 29623
$ws.Range("J55").Value = 29623
$ws.Range("L55").Value = 29623
$ws.Range("N55").Value = -30169

$ws.Range("H108").Value = 91995.8
$ws.Range("J108").Value = 91995.8
$ws.Range("L108").Value = 91995.8
$ws.Range("N108").Value = -99675.8

$ws.Range("H110").Value = 36679.75
$ws.Range("J110").Value = 36679.75
$ws.Range("L110").Value = 36679.75
$ws.Range("N110").Value = -44859.75

$ws.Range("H122").Value = 60369.375
$ws.Range("J122").Value = 60369.375
$ws.Range("L122").Value = 60369.375
$ws.Range("N122").Value = -70169.375

$ws.Range("H132").Value = 26464.072
$ws.Range("J132").Value = 26464.072
$ws.Range("L132").Value = 26464.072
$ws.Range("N132").Value = -36584.072

$ws.Range("H135").Value = 28181.727
$ws.Range("J135").Value = 28181.727
$ws.Range("L135").Value = 28181.727
$ws.Range("N135").Value = -38321.727

$ws.Range("H138").Value = 79824.5
$ws.Range("J138").Value = 79824.5
$ws.Range("L138").Value = 79824.5
$ws.Range("N138").Value = -90104.5

$ws.Range("H140").Value = 93492
$ws.Range("J140").Value = 93492
$ws.Range("L140").Value = 93492
$ws.Range("N140").Value = -103852

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 10113.6
$ws.Range("I7").Value = 12562.125
$ws.Range("J7").Value = 8481.25
$ws.Range("K7").Value = 12562.125
$ws.Range("L7").Value = 8481.25
$ws.Range("M7").Value = -12449.125
$ws.Range("N7").Value = -8707.25

$ws.Range("H9").Value = 79995.8
$ws.Range("J9").Value = 79995.8
$ws.Range("L9").Value = 79995.8
$ws.Range("N9").Value = -80331.8

$ws.Range("H18").Value = 24994.25
$ws.Range("J18").Value = 24994.25
$ws.Range("L18").Value = 24994.25
$ws.Range("N18").Value = -25454.25

$ws.Range("H70").Value = 23122.5
$ws.Range("J70").Value = 23122.5
$ws.Range("L70").Value = 23122.5
$ws.Range("N70").Value = -23752.5

$ws.Range("H73").Value = 23122.5
$ws.Range("J73").Value = 23122.5
$ws.Range("L73").Value = 23122.5
$ws.Range("N73").Value = -25306.5

$ws.Range("H114").Value = 51746
$ws.Range("J114").Value = 51746
$ws.Range("L114").Value = 51746
$ws.Range("N114").Value = -60424

$ws.Range("H117").Value = 33248.918
$ws.Range("J117").Value = 33248.918
$ws.Range("L117").Value = 33248.918
$ws.Range("N117").Value = -42426.918

$ws.Range("H118").Value = 64996
$ws.Range("J118").Value = 64996
$ws.Range("L118").Value = 64996
$ws.Range("N118").Value = -68310

$ws.Range("H132").Value = 2017263.4
$ws.Range("I132").Value = 2067709.8
$ws.Range("K132").Value = 6203129.4
$ws.Range("M132").Value = -6200599.4

$ws.Range("H134").Value = 2071701.2
$ws.Range("I134").Value = 2647725.8
$ws.Range("J134").Value = 127619
$ws.Range("K134").Value = 7943177.399999999
$ws.Range("L134").Value = 382857
$ws.Range("M134").Value = -7940642.399999999
$ws.Range("N134").Value = -387927

$ws.Range("H138").Value = 53918.4
$ws.Range("J138").Value = 54898
$ws.Range("L138").Value = 54898
$ws.Range("N138").Value = -65178

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 1683849.4
$ws.Range("J122").Value = 2525524
$ws.Range("L122").Value = 22729716
$ws.Range("N122").Value = -22734616

$ws.Range("H132").Value = 6440.1763
$ws.Range("J132").Value = 6798.9375
$ws.Range("L132").Value = 61190.4375
$ws.Range("N132").Value = -66250.4375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 17278.363
$ws.Range("J26").Value = 17278.363
$ws.Range("L26").Value = 17278.363
$ws.Range("N26").Value = -17838.363

$ws.Range("H50").Value = 17278.363
$ws.Range("J50").Value = 17278.363
$ws.Range("L50").Value = 17278.363
$ws.Range("N50").Value = -18274.363

$ws.Range("H59").Value = 6458.3335
$ws.Range("J59").Value = 7400
$ws.Range("L59").Value = 7400
$ws.Range("N59").Value = -8566

$ws.Range("H107").Value = 447.82608
$ws.Range("I107").Value = 391.4375
$ws.Range("J107").Value = 576.7143
$ws.Range("K107").Value = 391.4375
$ws.Range("L107").Value = 576.7143
$ws.Range("M107").Value = 1528.5625
$ws.Range("N107").Value = -4416.7143

$ws.Range("H108").Value = 51241.5
$ws.Range("J108").Value = 51241.5
$ws.Range("L108").Value = 51241.5
$ws.Range("N108").Value = -58921.5

$ws.Range("H109").Value = 36868.355
$ws.Range("I109").Value = 7900
$ws.Range("J109").Value = 39096.69
$ws.Range("K109").Value = 7900
$ws.Range("L109").Value = 39096.69
$ws.Range("M109").Value = -6860
$ws.Range("N109").Value = -41176.69

$ws.Range("H116").Value = 56766.668
$ws.Range("J116").Value = 58889
$ws.Range("L116").Value = 58889
$ws.Range("N116").Value = -68067

$ws.Range("H123").Value = 46499.5
$ws.Range("J123").Value = 46499.5
$ws.Range("L123").Value = 46499.5
$ws.Range("N123").Value = -51399.5

$ws.Range("H135").Value = 25000
$ws.Range("J135").Value = 25000
$ws.Range("L135").Value = 25000
$ws.Range("N135").Value = -35140

$ws.Range("H140").Value = 94332
$ws.Range("J140").Value = 94970.664
$ws.Range("L140").Value = 94970.664
$ws.Range("N140").Value = -105330.664

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H118").Value = 84011
$ws.Range("J118").Value = 84011
$ws.Range("L118").Value = 84011
$ws.Range("N118").Value = -87325

$ws.Range("H121").Value = 65430
$ws.Range("J121").Value = 65430
$ws.Range("L121").Value = 65430
$ws.Range("N121").Value = -68924

$ws.Range("H132").Value = 6520
$ws.Range("I132").Value = 6520
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 19560
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -17030
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 43000
$ws.Range("I75").Value = 25000
$ws.Range("K75").Value = 25000
$ws.Range("M75").Value = -24064

$ws.Range("H78").Value = 43000
$ws.Range("I78").Value = 25000
$ws.Range("K78").Value = 75000
$ws.Range("M78").Value = -70320

$ws.Range("H122").Value = 2067.111
$ws.Range("I122").Value = 648.75
$ws.Range("K122").Value = 1946.25
$ws.Range("M122").Value = 503.75

$ws.Range("H126").Value = 4528.067
$ws.Range("I126").Value = 4179.3125
$ws.Range("K126").Value = 12537.9375
$ws.Range("M126").Value = -10067.9375

$ws.Range("H136").Value = 1504
$ws.Range("I136").Value = 1504
$ws.Range("K136").Value = 4512
$ws.Range("M136").Value = -1962
